# Update NATMI ligand-receptor TPM statistics (Bmp2-Rgmb) for rows 2-37, columns E:T
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.2204535
$ws.Range("H2").Value = 8.440906999999999
$ws.Range("I2").Value = 0.3815978364461761
$ws.Range("J2").Value = 0.3430629039040183
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.8066855
$ws.Range("N2").Value = 25.613371
$ws.Range("O2").Value = 0.3362022094228136
$ws.Range("P2").Value = 0.2862725115007782
$ws.Range("Q2").Value = 54.05002064187425
$ws.Range("R2").Value = 216.200082567497
$ws.Range("S2").Value = 0.1282940357241699
$ws.Range("T2").Value = 0.09820947910335345

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.2204535
$ws.Range("H3").Value = 8.440906999999999
$ws.Range("I3").Value = 0.3815978364461761
$ws.Range("J3").Value = 0.3430629039040183
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.556209333333332
$ws.Range("N3").Value = 16.668628
$ws.Range("O3").Value = 0.1458620853836322
$ws.Range("P3").Value = 0.186299960314954
$ws.Range("Q3").Value = 23.44972312759933
$ws.Range("R3").Value = 140.698338765596
$ws.Range("S3").Value = 0.05566065620192145
$ws.Range("T3").Value = 0.0639126053828515

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.2204535
$ws.Range("H4").Value = 8.440906999999999
$ws.Range("I4").Value = 0.3815978364461761
$ws.Range("J4").Value = 0.3430629039040183
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.618423666666667
$ws.Range("N4").Value = 10.855271
$ws.Range("O4").Value = 0.0949911693670569
$ws.Range("P4").Value = 0.121325915756718
$ws.Range("Q4").Value = 15.27138882846617
$ws.Range("R4").Value = 91.628332970797
$ws.Range("S4").Value = 0.03624842471196119
$ws.Range("T4").Value = 0.04162242097831398

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.2204535
$ws.Range("H5").Value = 8.440906999999999
$ws.Range("I5").Value = 0.3815978364461761
$ws.Range("J5").Value = 0.3430629039040183
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 11.9979585
$ws.Range("N5").Value = 23.995917
$ws.Range("O5").Value = 0.3149714386492294
$ws.Range("P5").Value = 0.2681947419320252
$ws.Range("Q5").Value = 50.63682594417974
$ws.Range("R5").Value = 202.547303776719
$ws.Range("S5").Value = 0.1201924195308854
$ws.Range("T5").Value = 0.09200766697898938

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.2204535
$ws.Range("H6").Value = 8.440906999999999
$ws.Range("I6").Value = 0.3815978364461761
$ws.Range("J6").Value = 0.3430629039040183
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.6145123333333333
$ws.Range("N6").Value = 1.843537
$ws.Range("O6").Value = 0.01613223063721173
$ws.Range("P6").Value = 0.02060462744379138
$ws.Range("Q6").Value = 2.593520728009833
$ws.Range("R6").Value = 15.561124368059
$ws.Range("S6").Value = 0.006156024308210711
$ws.Range("T6").Value = 0.007068683324727499

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.2204535
$ws.Range("H7").Value = 8.440906999999999
$ws.Range("I7").Value = 0.3815978364461761
$ws.Range("J7").Value = 0.3430629039040183
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.498421666666667
$ws.Range("N7").Value = 10.495265
$ws.Range("O7").Value = 0.09184086654005638
$ws.Range("P7").Value = 0.1173022430517332
$ws.Range("Q7").Value = 14.76492596755917
$ws.Range("R7").Value = 88.58955580535499
$ws.Range("S7").Value = 0.03504627596902752
$ws.Range("T7").Value = 0.04024204813578255

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.087706
$ws.Range("H8").Value = 6.263118
$ws.Range("I8").Value = 0.1887626750859121
$ws.Range("J8").Value = 0.254551252439285
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.8066855
$ws.Range("N8").Value = 25.613371
$ws.Range("O8").Value = 0.3362022094228136
$ws.Range("P8").Value = 0.2862725115007782
$ws.Range("Q8").Value = 26.736594158463
$ws.Range("R8").Value = 160.419564950778
$ws.Range("S8").Value = 0.06346242842044433
$ws.Range("T8").Value = 0.07287102634146271

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.087706
$ws.Range("H9").Value = 6.263118
$ws.Range("I9").Value = 0.1887626750859121
$ws.Range("J9").Value = 0.254551252439285
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.556209333333332
$ws.Range("N9").Value = 16.668628
$ws.Range("O9").Value = 0.1458620853836322
$ws.Range("P9").Value = 0.186299960314954
$ws.Range("Q9").Value = 11.599731562456
$ws.Range("R9").Value = 104.397584062104
$ws.Range("S9").Value = 0.02753331743062414
$ws.Range("T9").Value = 0.04742288822756063

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.087706
$ws.Range("H10").Value = 6.263118
$ws.Range("I10").Value = 0.1887626750859121
$ws.Range("J10").Value = 0.254551252439285
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.618423666666667
$ws.Range("N10").Value = 10.855271
$ws.Range("O10").Value = 0.0949911693670569
$ws.Range("P10").Value = 0.121325915756718
$ws.Range("Q10").Value = 7.554204799442002
$ws.Range("R10").Value = 67.98784319497801
$ws.Range("S10").Value = 0.01793078723926461
$ws.Range("T10").Value = 0.03088366380921576

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.087706
$ws.Range("H11").Value = 6.263118
$ws.Range("I11").Value = 0.1887626750859121
$ws.Range("J11").Value = 0.254551252439285
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 11.9979585
$ws.Range("N11").Value = 23.995917
$ws.Range("O11").Value = 0.3149714386492294
$ws.Range("P11").Value = 0.2681947419320252
$ws.Range("Q11").Value = 25.048209948201
$ws.Range("R11").Value = 150.289259689206
$ws.Range("S11").Value = 0.05945485133508679
$ws.Range("T11").Value = 0.06826930745642786

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.087706
$ws.Range("H12").Value = 6.263118
$ws.Range("I12").Value = 0.1887626750859121
$ws.Range("J12").Value = 0.254551252439285
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.6145123333333333
$ws.Range("N12").Value = 1.843537
$ws.Range("O12").Value = 0.01613223063721173
$ws.Range("P12").Value = 0.02060462744379138
$ws.Range("Q12").Value = 1.282921085374
$ws.Range("R12").Value = 11.546289768366
$ws.Range("S12").Value = 0.003045163010182994
$ws.Range("T12").Value = 0.005244933721861958

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.087706
$ws.Range("H13").Value = 6.263118
$ws.Range("I13").Value = 0.1887626750859121
$ws.Range("J13").Value = 0.254551252439285
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.498421666666667
$ws.Range("N13").Value = 10.495265
$ws.Range("O13").Value = 0.09184086654005638
$ws.Range("P13").Value = 0.1173022430517332
$ws.Range("Q13").Value = 7.303675904030001
$ws.Range("R13").Value = 65.73308313627
$ws.Range("S13").Value = 0.01733612765030928
$ws.Range("T13").Value = 0.02985943288275611

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.100985
$ws.Range("H14").Value = 0.302955
$ws.Range("I14").Value = 0.009130691171817694
$ws.Range("J14").Value = 0.01231296850590131
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 12.8066855
$ws.Range("N14").Value = 25.613371
$ws.Range("O14").Value = 0.3362022094228136
$ws.Range("P14").Value = 0.2862725115007782
$ws.Range("Q14").Value = 1.2932831352175
$ws.Range("R14").Value = 7.759698811304999
$ws.Range("S14").Value = 0.003069758545522488
$ws.Range("T14").Value = 0.003524864418214352

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.100985
$ws.Range("H15").Value = 0.302955
$ws.Range("I15").Value = 0.009130691171817694
$ws.Range("J15").Value = 0.01231296850590131
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 5.556209333333332
$ws.Range("N15").Value = 16.668628
$ws.Range("O15").Value = 0.1458620853836322
$ws.Range("P15").Value = 0.186299960314954
$ws.Range("Q15").Value = 0.5610937995266665
$ws.Range("R15").Value = 5.049844195739999
$ws.Range("S15").Value = 0.001331821655315249
$ws.Range("T15").Value = 0.002293905544008692

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.100985
$ws.Range("H16").Value = 0.302955
$ws.Range("I16").Value = 0.009130691171817694
$ws.Range("J16").Value = 0.01231296850590131
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.618423666666667
$ws.Range("N16").Value = 10.855271
$ws.Range("O16").Value = 0.0949911693670569
$ws.Range("P16").Value = 0.121325915756718
$ws.Range("Q16").Value = 0.3654065139783333
$ws.Range("R16").Value = 3.288658625805
$ws.Range("S16").Value = 0.0008673350315404258
$ws.Range("T16").Value = 0.001493882179662104

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.100985
$ws.Range("H17").Value = 0.302955
$ws.Range("I17").Value = 0.009130691171817694
$ws.Range("J17").Value = 0.01231296850590131
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 11.9979585
$ws.Range("N17").Value = 23.995917
$ws.Range("O17").Value = 0.3149714386492294
$ws.Range("P17").Value = 0.2681947419320252
$ws.Range("Q17").Value = 1.2116138391225
$ws.Range("R17").Value = 7.269683034734999
$ws.Range("S17").Value = 0.002875906934249238
$ws.Range("T17").Value = 0.003302273410857355

# Row 18
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.100985
$ws.Range("H18").Value = 0.302955
$ws.Range("I18").Value = 0.009130691171817694
$ws.Range("J18").Value = 0.01231296850590131
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.6145123333333333
$ws.Range("N18").Value = 1.843537
$ws.Range("O18").Value = 0.01613223063721173
$ws.Range("P18").Value = 0.02060462744379138
$ws.Range("Q18").Value = 0.06205652798166666
$ws.Range("R18").Value = 0.5585087518349999
$ws.Range("S18").Value = 0.000147298415860916
$ws.Range("T18").Value = 0.0002537041287912329

# Row 19
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.100985
$ws.Range("H19").Value = 0.302955
$ws.Range("I19").Value = 0.009130691171817694
$ws.Range("J19").Value = 0.01231296850590131
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 3.498421666666667
$ws.Range("N19").Value = 10.495265
$ws.Range("O19").Value = 0.09184086654005638
$ws.Range("P19").Value = 0.1173022430517332
$ws.Range("Q19").Value = 0.3532881120083333
$ws.Range("R19").Value = 3.179593008075
$ws.Range("S19").Value = 0.0008385705893293799
$ws.Range("T19").Value = 0.001444338824367571

# Row 20
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 4.3548545
$ws.Range("H20").Value = 8.709709
$ws.Range("I20").Value = 0.393749879068184
$ws.Range("J20").Value = 0.3539877955886688
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 12.8066855
$ws.Range("N20").Value = 25.613371
$ws.Range("O20").Value = 0.3362022094228136
$ws.Range("P20").Value = 0.2862725115007782
$ws.Range("Q20").Value = 55.77125197975975
$ws.Range("R20").Value = 223.085007919039
$ws.Range("S20").Value = 0.1323795793026891
$ws.Range("T20").Value = 0.1013369752837923

# Row 21
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 4.3548545
$ws.Range("H21").Value = 8.709709
$ws.Range("I21").Value = 0.393749879068184
$ws.Range("J21").Value = 0.3539877955886688
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 5.556209333333332
$ws.Range("N21").Value = 16.668628
$ws.Range("O21").Value = 0.1458620853836322
$ws.Range("P21").Value = 0.186299960314954
$ws.Range("Q21").Value = 24.19648321820866
$ws.Range("R21").Value = 145.178899309252
$ws.Range("S21").Value = 0.05743317848043831
$ws.Range("T21").Value = 0.06594791227014705

# Row 22
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 4.3548545
$ws.Range("H22").Value = 8.709709
$ws.Range("I22").Value = 0.393749879068184
$ws.Range("J22").Value = 0.3539877955886688
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 3.618423666666667
$ws.Range("N22").Value = 10.855271
$ws.Range("O22").Value = 0.0949911693670569
$ws.Range("P22").Value = 0.121325915756718
$ws.Range("Q22").Value = 15.75770858768984
$ws.Range("R22").Value = 94.54625152613902
$ws.Range("S22").Value = 0.03740276145082404
$ws.Range("T22").Value = 0.04294789346649715

# Row 23
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 4.3548545
$ws.Range("H23").Value = 8.709709
$ws.Range("I23").Value = 0.393749879068184
$ws.Range("J23").Value = 0.3539877955886688
$ws.Range("K23").Value = 2
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 11.9979585
$ws.Range("N23").Value = 23.995917
$ws.Range("O23").Value = 0.3149714386492294
$ws.Range("P23").Value = 0.2681947419320252
$ws.Range("Q23").Value = 52.24936356453825
$ws.Range("R23").Value = 208.997454258153
$ws.Range("S23").Value = 0.124019965878066
$ws.Range("T23").Value = 0.09493766548498954

# Row 24
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 4.3548545
$ws.Range("H24").Value = 8.709709
$ws.Range("I24").Value = 0.393749879068184
$ws.Range("J24").Value = 0.3539877955886688
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.6145123333333333
$ws.Range("N24").Value = 1.843537
$ws.Range("O24").Value = 0.01613223063721173
$ws.Range("P24").Value = 0.02060462744379138
$ws.Range("Q24").Value = 2.676111800122167
$ws.Range("R24").Value = 16.056670800733
$ws.Range("S24").Value = 0.006352063862502169
$ws.Range("T24").Value = 0.007293786647753498

# Row 25
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 4.3548545
$ws.Range("H25").Value = 8.709709
$ws.Range("I25").Value = 0.393749879068184
$ws.Range("J25").Value = 0.3539877955886688
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 3.498421666666667
$ws.Range("N25").Value = 10.495265
$ws.Range("O25").Value = 0.09184086654005638
$ws.Range("P25").Value = 0.1173022430517332
$ws.Range("Q25").Value = 15.23511733798083
$ws.Range("R25").Value = 91.410704027885
$ws.Range("S25").Value = 0.03616233009366442
$ws.Range("T25").Value = 0.04152356243548928

# Row 26
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 0.3333333333333333
$ws.Range("G26").Value = 0.06028266666666667
$ws.Range("H26").Value = 0.180848
$ws.Range("I26").Value = 0.005450536340515544
$ws.Range("J26").Value = 0.007350186424898878
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 12.8066855
$ws.Range("N26").Value = 25.613371
$ws.Range("O26").Value = 0.3362022094228136
$ws.Range("P26").Value = 0.2862725115007782
$ws.Range("Q26").Value = 0.7720211531013335
$ws.Range("R26").Value = 4.632126918608001
$ws.Range("S26").Value = 0.001832482360220663
$ws.Range("T26").Value = 0.002104156327854728

# Row 27
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 0.3333333333333333
$ws.Range("G27").Value = 0.06028266666666667
$ws.Range("H27").Value = 0.180848
$ws.Range("I27").Value = 0.005450536340515544
$ws.Range("J27").Value = 0.007350186424898878
$ws.Range("K27").Value = 3
$ws.Range("L27").Value = 1
$ws.Range("M27").Value = 5.556209333333332
$ws.Range("N27").Value = 16.668628
$ws.Range("O27").Value = 0.1458620853836322
$ws.Range("P27").Value = 0.186299960314954
$ws.Range("Q27").Value = 0.3349431151715555
$ws.Range("R27").Value = 3.014488036544
$ws.Range("S27").Value = 0.0007950265970868686
$ws.Range("T27").Value = 0.001369339439266175

# Row 28
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = 0.3333333333333333
$ws.Range("G28").Value = 0.06028266666666667
$ws.Range("H28").Value = 0.180848
$ws.Range("I28").Value = 0.005450536340515544
$ws.Range("J28").Value = 0.007350186424898878
$ws.Range("K28").Value = 3
$ws.Range("L28").Value = 1
$ws.Range("M28").Value = 3.618423666666667
$ws.Range("N28").Value = 10.855271
$ws.Range("O28").Value = 0.0949911693670569
$ws.Range("P28").Value = 0.121325915756718
$ws.Range("Q28").Value = 0.2181282277564445
$ws.Range("R28").Value = 1.963154049808
$ws.Range("S28").Value = 0.0005177528206632106
$ws.Range("T28").Value = 0.0008917680989834537

# Row 29
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 0.3333333333333333
$ws.Range("G29").Value = 0.06028266666666667
$ws.Range("H29").Value = 0.180848
$ws.Range("I29").Value = 0.005450536340515544
$ws.Range("J29").Value = 0.007350186424898878
$ws.Range("K29").Value = 2
$ws.Range("L29").Value = 1
$ws.Range("M29").Value = 11.9979585
$ws.Range("N29").Value = 23.995917
$ws.Range("O29").Value = 0.3149714386492294
$ws.Range("P29").Value = 0.2681947419320252
$ws.Range("Q29").Value = 0.723268932936
$ws.Range("R29").Value = 4.339613597616
$ws.Range("S29").Value = 0.001716763272582087
$ws.Range("T29").Value = 0.00197128135137803

# Row 30
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = 0.3333333333333333
$ws.Range("G30").Value = 0.06028266666666667
$ws.Range("H30").Value = 0.180848
$ws.Range("I30").Value = 0.005450536340515544
$ws.Range("J30").Value = 0.007350186424898878
$ws.Range("K30").Value = 2
$ws.Range("L30").Value = 0.6666666666666666
$ws.Range("M30").Value = 0.6145123333333333
$ws.Range("N30").Value = 1.843537
$ws.Range("O30").Value = 0.01613223063721173
$ws.Range("P30").Value = 0.02060462744379138
$ws.Range("Q30").Value = 0.03704444215288889
$ws.Range("R30").Value = 0.333399979376
$ws.Range("S30").Value = [double]"8.792930934170075E-05"
$ws.Range("T30").Value = 0.0001514478529274543

# Row 31
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 0.3333333333333333
$ws.Range("G31").Value = 0.06028266666666667
$ws.Range("H31").Value = 0.180848
$ws.Range("I31").Value = 0.005450536340515544
$ws.Range("J31").Value = 0.007350186424898878
$ws.Range("K31").Value = 3
$ws.Range("L31").Value = 1
$ws.Range("M31").Value = 3.498421666666667
$ws.Range("N31").Value = 10.495265
$ws.Range("O31").Value = 0.09184086654005638
$ws.Range("P31").Value = 0.1173022430517332
$ws.Range("Q31").Value = 0.2108941871911111
$ws.Range("R31").Value = 1.89804768472
$ws.Range("S31").Value = 0.0005005819806210154
$ws.Range("T31").Value = 0.0008621933544890382

# Row 32
$ws.Range("E32").Value = 2
$ws.Range("F32").Value = 0.6666666666666666
$ws.Range("G32").Value = 0.2356696666666667
$ws.Range("H32").Value = 0.707009
$ws.Range("I32").Value = 0.02130838188739468
$ws.Range("J32").Value = 0.02873489313722756
$ws.Range("K32").Value = 2
$ws.Range("L32").Value = 1
$ws.Range("M32").Value = 12.8066855
$ws.Range("N32").Value = 25.613371
$ws.Range("O32").Value = 0.3362022094228136
$ws.Range("P32").Value = 0.2862725115007782
$ws.Range("Q32").Value = 3.018147302889834
$ws.Range("R32").Value = 18.108883817339
$ws.Range("S32").Value = 0.007163925069767155
$ws.Range("T32").Value = 0.00822601002610061

# Row 33
$ws.Range("E33").Value = 2
$ws.Range("F33").Value = 0.6666666666666666
$ws.Range("G33").Value = 0.2356696666666667
$ws.Range("H33").Value = 0.707009
$ws.Range("I33").Value = 0.02130838188739468
$ws.Range("J33").Value = 0.02873489313722756
$ws.Range("K33").Value = 3
$ws.Range("L33").Value = 1
$ws.Range("M33").Value = 5.556209333333332
$ws.Range("N33").Value = 16.668628
$ws.Range("O33").Value = 0.1458620853836322
$ws.Range("P33").Value = 0.186299960314954
$ws.Range("Q33").Value = 1.309430001516889
$ws.Range("R33").Value = 11.784870013652
$ws.Range("S33").Value = 0.003108085018246205
$ws.Range("T33").Value = 0.005353309451119939

# Row 34
$ws.Range("E34").Value = 2
$ws.Range("F34").Value = 0.6666666666666666
$ws.Range("G34").Value = 0.2356696666666667
$ws.Range("H34").Value = 0.707009
$ws.Range("I34").Value = 0.02130838188739468
$ws.Range("J34").Value = 0.02873489313722756
$ws.Range("K34").Value = 3
$ws.Range("L34").Value = 1
$ws.Range("M34").Value = 3.618423666666667
$ws.Range("N34").Value = 10.855271
$ws.Range("O34").Value = 0.0949911693670569
$ws.Range("P34").Value = 0.121325915756718
$ws.Range("Q34").Value = 0.8527526993821112
$ws.Range("R34").Value = 7.674774294439001
$ws.Range("S34").Value = 0.002024108112803436
$ws.Range("T34").Value = 0.003486287224045566

# Row 35
$ws.Range("E35").Value = 2
$ws.Range("F35").Value = 0.6666666666666666
$ws.Range("G35").Value = 0.2356696666666667
$ws.Range("H35").Value = 0.707009
$ws.Range("I35").Value = 0.02130838188739468
$ws.Range("J35").Value = 0.02873489313722756
$ws.Range("K35").Value = 2
$ws.Range("L35").Value = 1
$ws.Range("M35").Value = 11.9979585
$ws.Range("N35").Value = 23.995917
$ws.Range("O35").Value = 0.3149714386492294
$ws.Range("P35").Value = 0.2681947419320252
$ws.Range("Q35").Value = 2.8275548803755
$ws.Range("R35").Value = 16.965329282253
$ws.Range("S35").Value = 0.006711531698359886
$ws.Range("T35").Value = 0.00770654724938307

# Row 36
$ws.Range("E36").Value = 2
$ws.Range("F36").Value = 0.6666666666666666
$ws.Range("G36").Value = 0.2356696666666667
$ws.Range("H36").Value = 0.707009
$ws.Range("I36").Value = 0.02130838188739468
$ws.Range("J36").Value = 0.02873489313722756
$ws.Range("K36").Value = 2
$ws.Range("L36").Value = 0.6666666666666666
$ws.Range("M36").Value = 0.6145123333333333
$ws.Range("N36").Value = 1.843537
$ws.Range("O36").Value = 0.01613223063721173
$ws.Range("P36").Value = 0.02060462744379138
$ws.Range("Q36").Value = 0.1448219167592222
$ws.Range("R36").Value = 1.303397250833
$ws.Range("S36").Value = 0.0003437517311132359
$ws.Range("T36").Value = 0.0005920717677297315

# Row 37
$ws.Range("E37").Value = 2
$ws.Range("F37").Value = 0.6666666666666666
$ws.Range("G37").Value = 0.2356696666666667
$ws.Range("H37").Value = 0.707009
$ws.Range("I37").Value = 0.02130838188739468
$ws.Range("J37").Value = 0.02873489313722756
$ws.Range("K37").Value = 3
$ws.Range("L37").Value = 1
$ws.Range("M37").Value = 3.498421666666667
$ws.Range("N37").Value = 10.495265
$ws.Range("O37").Value = 0.09184086654005638
$ws.Range("P37").Value = 0.1173022430517332
$ws.Range("Q37").Value = 0.8244718680427778
$ws.Range("R37").Value = 7.420246812385
$ws.Range("S37").Value = 0.00195698025710477
$ws.Range("T37").Value = 0.003370667418848648

